$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update instruction text (row 4 & 5 swap content; rows 6 & 7 get wording tweaks) ---
$ws.Range("A4").Value = "Some images are special: the correct response is the opposite of the preceding correct response. That is, if the previous correct answer was “Left” then the correct response would now be “Right”."
$ws.Range("A5").Value = "After each response, you will be told whether you got the item correct or incorrect."
$ws.Range("A6").Value = "You can take as long as you like on each image, but the task will not continue until you press the “Left” or “Right” arrow key."
$ws.Range("A7").Value = "The experiment will have three blocks, each block will take approximately 10 minutes to complete."

# --- Reshape the trailing blank rows: remove row 8, add blank styled rows 9 and 11 ---
$ws.Rows(8).Delete()

$ws.Range("A9").Value = " "
$ws.Range("A9").ClearContents()
$ws.Range("A9").VerticalAlignment = -4108

$ws.Range("A11").Value = " "
$ws.Range("A11").ClearContents()
$ws.Range("A11").VerticalAlignment = -4108
